$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "Anna Zandonati" team as row 64 (sheet was A1:F63, now A1:F64)
$ws.Range("A64").Value = "Anna Zandonati"
$ws.Range("B64").Value = "Stefano Tita | Clitoriders"
$ws.Range("C64").Value = "Mattia Festi | Shark Attack"
$ws.Range("D64").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E64").Value = "Federico  Manica | iMontagna"
$ws.Range("F64").Value = "Alessandro  Tengattini | Herta Vernello"
